# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet gains a new "property_category" column
# (always "stock") inserted right after the "total" column and before the
# existing "date" column. All columns from the old "date" column onward
# (date, legislator_name, legislator_id) shift one slot to the right.
# Also fixes a stray space in one of the company names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new blank column at H; everything from H onward (date,
# legislator_name, legislator_id) slides right to I, J, K and keeps its
# original formatting/style.
$ws.Columns("H").Insert()

# New header + values for the inserted "property_category" column.
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(2, 8).Value = "stock"
$ws.Cells.Item(3, 8).Value = "stock"
$ws.Cells.Item(4, 8).Value = "stock"

# Fix typo: remove stray space in company name.
$ws.Cells.Item(2, 2).Value = "杏輝藥品工業股份有限公司"
